# Update data for June 2024 baseline
# Applies revised EMTR figures (columns B:L) to the data rows of the
# "EMTR on Capital Income" sheet, removes the now-unused footnote row
# and its associated conditional formatting / shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for each affected data row, columns B (2024) through L (2034).
$emtrUpdates = @{
    12 = @(14.46, 14.86, 14.49, 15.07, 15.09, 15.05, 15.15, 15.19, 15.24, 15.25, 15.29)
    14 = @(17.82, 18.53, 21.03, 21.93, 21.96, 21.98, 22.17, 22.21, 22.25, 22.27, 22.32)
    16 = @(10.21, 12.49, 15.31, 18.13, 18.25, 18.37, 18.74, 18.81, 18.88, 18.89, 18.96)
    17 = @(19.51, 19.9,  21.81, 22.4,  22.37, 22.32, 22.43, 22.45, 22.45, 22.5,  22.53)
    18 = @(23.96, 24.01, 27.88, 27.91, 27.86, 27.78, 27.74, 27.74, 27.78, 27.79, 27.82)
    19 = @(-7.97, -7.72, -6.37, -6.06, -5.54, -4.47, -2.71, -2.52, -2.41, -2.4,  -2.31)
    20 = @(17.19, 20.66, 24.39, 27.34, 27.4,  27.4,  27.59, 27.65, 27.71, 27.73, 27.79)
    21 = @(29.5,  29.48, 31.82, 31.79, 31.77, 31.69, 31.71, 31.73, 31.78, 31.79, 31.83)
    23 = @(20.54, 21.25, 24.11, 24.98, 24.96, 24.99, 25.18, 25.19, 25.21, 25.24, 25.26)
    24 = @(5.45,  6.15,  6.15,  7.27,  7.53,  7.46,  7.74,  7.91,  8.02,  8.02,  8.18)
    26 = @(16.63, 17.39, 18.58, 19.59, 19.66, 19.72, 20.04, 20.1,  20.13, 20.17, 20.22)
    27 = @(20.78, 21.37, 26.8,  27.47, 27.4,  27.33, 27.27, 27.26, 27.3,  27.3,  27.32)
    29 = @(8.28,  8.04,  0.96,  0.68,  0.67,  0.51,  0.37,  0.4,   0.48,  0.43,  0.45)
    31 = @(-0.15, -0.15, -2.23, -2.26, -2.26, -2.26, -2.26, -2.26, -2.26, -2.27, -2.27)
    32 = @(26.85, 26.25, 9.72,  8.79,  8.79,  8.17,  7.67,  7.78,  8.09,  7.91,  8.01)
    35 = @(15.13, 15.38, 14.71, 15.09, 15.09, 15.03, 15.09, 15.12, 15.17, 15.17, 15.21)
    36 = @(19.85, 20.39, 23.09, 23.78, 23.78, 23.78, 23.93, 23.96, 24,    24.02, 24.06)
    37 = @(27.11, 27.1,  30.44, 30.42, 30.38, 30.3,  30.29, 30.3,  30.34, 30.35, 30.38)
    38 = @(8.28,  8.04,  0.96,  0.68,  0.67,  0.51,  0.37,  0.4,   0.48,  0.43,  0.45)
}

foreach ($rowNum in $emtrUpdates.Keys) {
    $values = $emtrUpdates[$rowNum]
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Column B is index 2; data spans B..L (columns 2..12)
        $col = $i + 2
        $ws.Cells.Item($rowNum, $col).Value = $values[$i]
    }
}

# Remove the footnote row (previously "* = between zero and 0.05 percentage
# points.") entirely -- it is no longer referenced now that the conditional
# formatting that used it has been removed.
$ws.Range("A43:L43").Clear()

# Remove the conditional formatting rule (highlighting near-zero EMTRs)
# that applied to the data block; its backing dxf is no longer needed.
$ws.Range("B12:L38").FormatConditions.Delete()
